$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing D:K data to E:L
$ws.Columns("D:D").Insert()

# Copy number formatting from column E (the shifted original D) into new column D
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)

# Populate the new column D with the newest reporting period's data
$ws.Range("D7").Value = 43463
$ws.Range("D8").Value = 398800
$ws.Range("D9").Value = 179400
$ws.Range("D10").Value = 219400
$ws.Range("D12").Value = 82200
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 31400
$ws.Range("D15").Value = 18000
$ws.Range("D17").Value = 401900
$ws.Range("D18").Value = -3100
$ws.Range("D20").Value = -200
$ws.Range("D21").Value = 35900
$ws.Range("D22").Value = 20600
$ws.Range("D23").Value = -24000
$ws.Range("D24").Value = 2400
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -26300
$ws.Range("D27").Value = -26300
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 200
$ws.Range("D33").Value = -26300
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -26300
$ws.Range("D38").Value = 43463
$ws.Range("D41").Value = 119100
$ws.Range("D42").Value = 9600
$ws.Range("D43").Value = 60900
$ws.Range("D44").Value = 67100
$ws.Range("D45").Value = 27800
$ws.Range("D46").Value = 284400
$ws.Range("D47").Value = "NA"
$ws.Range("D48").Value = 34900
$ws.Range("D49").Value = 288800
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 15500
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 623700
$ws.Range("D57").Value = 31900
$ws.Range("D58").Value = 8300
$ws.Range("D59").Value = 29200
$ws.Range("D60").Value = 69400
$ws.Range("D61").Value = 251400
$ws.Range("D62").Value = 44500
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 365200
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -476800
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 258500
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43463
$ws.Range("D81").Value = -26300
$ws.Range("D83").Value = 39300
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 51500
$ws.Range("D91").Value = -8400
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -21100
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -16800
$ws.Range("D101").Value = -1300
$ws.Range("D102").Value = 12200

$excel.CutCopyMode = 0
